$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E4").Value = 23
$ws.Range("E15").Value = 158
$ws.Range("F15").Value = 84
$ws.Range("H15").Value = 125
$ws.Range("F17").Value = 57
$ws.Range("H17").Value = 89
$ws.Range("F18").Value = 47
$ws.Range("H18").Value = 83
$ws.Range("F19").Value = 29
$ws.Range("H19").Value = 42
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = 7
$ws.Range("E25").Value = 23
$ws.Range("F26").Value = 13
$ws.Range("H26").Value = 23
$ws.Range("E28").Value = 17
$ws.Range("F29").Value = 10
$ws.Range("H29").Value = 13
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = 3
$ws.Range("H35").Value = 4
$ws.Range("F38").Value = 17
$ws.Range("H38").Value = 37
$ws.Range("F39").Value = 15
$ws.Range("H39").Value = 23
$ws.Range("E47").Value = 59
$ws.Range("F47").Value = 37
$ws.Range("H47").Value = 48
$ws.Range("F48").Value = 21
$ws.Range("H48").Value = 26
$ws.Range("F52").Value = 3
$ws.Range("H52").Value = 3
$ws.Range("F68").Value = 10
$ws.Range("H68").Value = 14
$ws.Range("F70").Value = 21
$ws.Range("H70").Value = 33
$ws.Range("E72").Value = 43
$ws.Range("F72").Value = 21
$ws.Range("H72").Value = 32
$ws.Range("E74").Value = 20
$ws.Range("F76").Value = 19
$ws.Range("H76").Value = 36
$ws.Range("E77").Value = 55
$ws.Range("F77").Value = 21
$ws.Range("H77").Value = 38
$ws.Range("E80").Value = 28
$ws.Range("E88").Value = 22
$ws.Range("F89").Value = 18
$ws.Range("H89").Value = 25
